# Generate Report for Handoff
#
# The localization status moves from "In Translation" to
# "Ready for handoff" and the corresponding "latest handoff" timestamps
# are refreshed. This touches the three report sheets (Overview, zh-cn,
# de-de). Because the new status text is longer than the old text, the
# status columns also widen to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Columns E (zh-cn) and F (de-de) hold the per-locale status, column G
# holds the latest handoff xliff generation timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-25 19:01:21"

# Widen the now-longer status columns to match the new text length.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet -------------------------------------------------------
# Column C is the Status column, column H is the Latest Handoff Datetime.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-25 19:01:16"
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ---------------------------------------------------------
# Column C is the Status column, column H is the Latest Handoff Datetime.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-25 19:01:21"
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
